$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3538.3333
$ws.Range("I62").Value = 2685
$ws.Range("J62").Value = 5245
$ws.Range("K62").Value = 2685
$ws.Range("L62").Value = 5245
$ws.Range("M62").Value = -2061
$ws.Range("N62").Value = -6493
$ws.Range("H65").Value = 3538.3333
$ws.Range("I65").Value = 2685
$ws.Range("J65").Value = 5245
$ws.Range("K65").Value = 13425
$ws.Range("L65").Value = 26225
$ws.Range("M65").Value = -10305
$ws.Range("N65").Value = -32465
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H100").Value = 33334638
$ws.Range("I100").Value = 33334638
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 33334638
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -33334097
$ws.Range("N100").ClearContents()
$ws.Range("H112").Value = 1261.5344
$ws.Range("J112").Value = 1308.5272
$ws.Range("L112").Value = 3925.5816
$ws.Range("N112").Value = -6141.5816
$ws.Range("H115").Value = 1257.5
$ws.Range("I115").Value = 1257.5
$ws.Range("K115").Value = 3772.5
$ws.Range("M115").Value = -2205.5
$ws.Range("H135").Value = 910.13336
$ws.Range("I135").Value = 504
$ws.Range("K135").Value = 4536
$ws.Range("M135").Value = -2001
$ws.Range("H137").Value = 3404184
$ws.Range("I137").Value = 4763507.5
$ws.Range("J137").Value = 5875
$ws.Range("K137").Value = 14290522.5
$ws.Range("L137").Value = 17625
$ws.Range("M137").Value = -14287972.5
$ws.Range("N137").Value = -22725
$ws.Range("H138").Value = 2644.54
$ws.Range("I138").Value = 863.375
$ws.Range("J138").Value = 2983.8096
$ws.Range("K138").Value = 2590.125
$ws.Range("L138").Value = 8951.4288
$ws.Range("M138").Value = 2549.875
$ws.Range("N138").Value = -19231.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 40998.4
$ws.Range("J47").Value = 40998.4
$ws.Range("L47").Value = 40998.4
$ws.Range("N47").Value = -42448.4
$ws.Range("H61").Value = 1041.9231
$ws.Range("I61").Value = 911.5
$ws.Range("J61").Value = 2607
$ws.Range("K61").Value = 911.5
$ws.Range("L61").Value = 2607
$ws.Range("M61").Value = -699.5
$ws.Range("N61").Value = -3031
$ws.Range("H74").Value = 7679.2666
$ws.Range("I74").Value = 11062
$ws.Range("J74").Value = 3813.2856
$ws.Range("K74").Value = 11062
$ws.Range("L74").Value = 3813.2856
$ws.Range("M74").Value = -10188
$ws.Range("N74").Value = -5561.2856
$ws.Range("H77").Value = 7679.2666
$ws.Range("I77").Value = 11062
$ws.Range("J77").Value = 3813.2856
$ws.Range("K77").Value = 55310
$ws.Range("L77").Value = 19066.428
$ws.Range("M77").Value = -50942
$ws.Range("N77").Value = -27802.428
$ws.Range("H88").Value = 7410457
$ws.Range("I88").Value = 13335863
$ws.Range("J88").Value = 3700
$ws.Range("K88").Value = 13335863
$ws.Range("L88").Value = 3700
$ws.Range("M88").Value = -13335457
$ws.Range("N88").Value = -4512
$ws.Range("H91").Value = 7410457
$ws.Range("I91").Value = 13335863
$ws.Range("J91").Value = 3700
$ws.Range("K91").Value = 13335863
$ws.Range("L91").Value = 3700
$ws.Range("M91").Value = -13334459
$ws.Range("N91").Value = -6508
$ws.Range("H110").Value = 701.1539
$ws.Range("I110").Value = 770.2
$ws.Range("J110").Value = 471
$ws.Range("K110").Value = 770.2
$ws.Range("L110").Value = 471
$ws.Range("M110").Value = 1274.8
$ws.Range("N110").Value = -4561
$ws.Range("H136").Value = 1041.9231
$ws.Range("I136").Value = 911.5
$ws.Range("J136").Value = 2607
$ws.Range("K136").Value = 2734.5
$ws.Range("L136").Value = 7821
$ws.Range("M136").Value = -184.5
$ws.Range("N136").Value = -12921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 28218.334
$ws.Range("I20").Value = 4905
$ws.Range("J20").Value = 39875
$ws.Range("K20").Value = 4905
$ws.Range("L20").Value = 39875
$ws.Range("M20").Value = -4658
$ws.Range("N20").Value = -40369
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H94").Value = 1224.1666
$ws.Range("I94").Value = 869
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 869
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -418
$ws.Range("N94").Value = -3902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5422.222
$ws.Range("I132").Value = 2937.3333
$ws.Range("J132").Value = 6664.6665
$ws.Range("K132").Value = 8811.999899999999
$ws.Range("L132").Value = 19993.9995
$ws.Range("M132").Value = -6281.999899999999
$ws.Range("N132").Value = -25053.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 797.8333
$ws.Range("I92").Value = 798
$ws.Range("K92").Value = 2394
$ws.Range("M92").Value = -1146
$ws.Range("H104").Value = 2139.8572
$ws.Range("I104").Value = 3000
$ws.Range("J104").Value = 1996.5
$ws.Range("K104").Value = 9000
$ws.Range("L104").Value = 5989.5
$ws.Range("M104").Value = -6379
$ws.Range("N104").Value = -11231.5
$ws.Range("H123").Value = 3603.9
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H131").Value = 706.05
$ws.Range("I131").Value = 278
$ws.Range("J131").Value = 813.0625
$ws.Range("K131").Value = 834
$ws.Range("L131").Value = 2439.1875
$ws.Range("M131").Value = 4206
$ws.Range("N131").Value = -12519.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6740.4
$ws.Range("I70").Value = 5913.75
$ws.Range("J70").Value = 8210
$ws.Range("K70").Value = 5913.75
$ws.Range("L70").Value = 8210
$ws.Range("M70").Value = -5643.75
$ws.Range("N70").Value = -8750
$ws.Range("H73").Value = 6740.4
$ws.Range("I73").Value = 5913.75
$ws.Range("J73").Value = 8210
$ws.Range("K73").Value = 5913.75
$ws.Range("L73").Value = 8210
$ws.Range("M73").Value = -4977.75
$ws.Range("N73").Value = -10082
$ws.Range("H80").Value = 83336000
$ws.Range("I80").Value = 125002500
$ws.Range("K80").Value = 125002500
$ws.Range("M80").Value = -125001502
$ws.Range("H83").Value = 83336000
$ws.Range("I83").Value = 125002500
$ws.Range("K83").Value = 625012500
$ws.Range("M83").Value = -625007508
$ws.Range("H97").Value = 1362
$ws.Range("I97").Value = 1234
$ws.Range("K97").Value = 1234
$ws.Range("M97").Value = -738
$ws.Range("H102").Value = 3783.2
$ws.Range("I102").Value = 2342
$ws.Range("J102").Value = 5224.4
$ws.Range("K102").Value = 2342
$ws.Range("L102").Value = 5224.4
$ws.Range("M102").Value = -720
$ws.Range("N102").Value = -8468.4
$ws.Range("H122").Value = 7656.3335
$ws.Range("I122").Value = 2381.4
$ws.Range("K122").Value = 7144.200000000001
$ws.Range("M122").Value = -4694.200000000001
$ws.Range("H132").Value = 2579.3823
$ws.Range("I132").Value = 1255.8572
$ws.Range("J132").Value = 4717.385
$ws.Range("K132").Value = 3767.5716
$ws.Range("L132").Value = 14152.155
$ws.Range("M132").Value = -1237.5716
$ws.Range("N132").Value = -19212.155
$ws.Range("H133").Value = 41586.668
$ws.Range("J133").Value = 41586.668
$ws.Range("L133").Value = 41586.668
$ws.Range("N133").Value = -51706.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 44998.5
$ws.Range("J47").Value = 44998.5
$ws.Range("L47").Value = 44998.5
$ws.Range("N47").Value = -45978.5
$ws.Range("H52").Value = 44998.5
$ws.Range("J52").Value = 44998.5
$ws.Range("L52").Value = 44998.5
$ws.Range("N52").Value = -45464.5
$ws.Range("H82").Value = 1236.6666
$ws.Range("I82").Value = 677.5625
$ws.Range("J82").Value = 2049.9092
$ws.Range("K82").Value = 677.5625
$ws.Range("L82").Value = 2049.9092
$ws.Range("M82").Value = -316.5625
$ws.Range("N82").Value = -2771.9092
$ws.Range("H85").Value = 1236.6666
$ws.Range("I85").Value = 677.5625
$ws.Range("J85").Value = 2049.9092
$ws.Range("K85").Value = 677.5625
$ws.Range("L85").Value = 2049.9092
$ws.Range("M85").Value = 570.4375
$ws.Range("N85").Value = -4545.9092
$ws.Range("H122").Value = 4063.75
$ws.Range("I122").Value = 1885
$ws.Range("K122").Value = 5655
$ws.Range("M122").Value = -3205

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13891216
$ws.Range("I132").Value = 1726.95
$ws.Range("J132").Value = 83338664
$ws.Range("K132").Value = 5180.85
$ws.Range("L132").Value = 250015992
$ws.Range("M132").Value = -2650.85
$ws.Range("N132").Value = -250021052

